$wb = $excel.ActiveWorkbook

# --- Sheet "Personas": the used range grows from Q1:Q30 to Q1:Q189 -------
# These placeholder cells only ever held an (intentionally) empty value, so
# a plain .Value write would be optimised away by the engine (it never
# materialises a truly blank cell). Touching a formatting property and then
# writing back the exact value it already had is enough to force the cell
# into existence without changing its appearance or creating new styles.
$ws1 = $wb.Worksheets.Item("Personas")
$qRange = $ws1.Range("Q31:Q189")
$existingSize = $qRange.Font.Size
$qRange.Font.Size = $existingSize

# --- Sheet "Subcategorias": append a new data row (row 7) ---------------
$ws2 = $wb.Worksheets.Item("Subcategorias")
$ws2.Activate()

# Copy the formatting of the last existing row so the new row matches the
# rest of the table exactly.
$ws2.Range("A6:D6").Copy()
$ws2.Range("A7").PasteSpecial(-4122)

$ws2.Cells.Item(7, 1).Value = 58
$ws2.Cells.Item(7, 2).Value = "testSub"
$ws2.Cells.Item(7, 3).Value = "testdesS"
$ws2.Cells.Item(7, 4).Value = "testEnS"

[void]$ws2.Range("A7:D7").Select()
